# Add columns I (I0) and J (IF) to the worksheet, mirroring the style of
# the existing header/data columns (H = IP).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (row 1) - copy style from H1 (bold, bordered, centered header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Data values for columns I and J, rows 2-16
$values = @(
    @(7, 7),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(7, 8),
    @(6, 7),
    @(8, 9),
    @(5, 5),
    @(7, 7),
    @(8, 8),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $pair = $values[$i]
    $ws.Cells.Item($row, 9).Value = $pair[0]   # column I
    $ws.Cells.Item($row, 10).Value = $pair[1]  # column J
}
